# B6-PowerPoint.pptx edit
#
# 1) Three tables (slides 14, 15, 16) switch from the custom table style
#    {909E6BE8-1D49-4FA2-A262-76BA387341D0} to the built-in table style
#    {EC4EFC09-0E0E-4E82-962C-551781366F3A}.
# 2) The deck's theme colour scheme (the "Integral"/Red Violet palette
#    driving every slide) is replaced with the stock "Office Theme"
#    palette (the deck's two theme parts effectively swap colour
#    content).

$p = $ppt.ActivePresentation

# --- 1) Re-style the three tables -----------------------------------
$tableSlides = 14,15,16
foreach ($n in $tableSlides) {
    $slide = $p.Slides.Item($n)
    $tbl = $slide.Shapes.Item(1).Table
    $tbl.ApplyStyle("{EC4EFC09-0E0E-4E82-962C-551781366F3A}")
}

# --- 2) Swap the theme colour scheme back to the stock Office colours -
$tcs = $p.Slides.Item(1).ThemeColorScheme
$tcs.Colors(1).RGB  = 0         # dk1      000000
$tcs.Colors(2).RGB  = 16777215  # lt1      FFFFFF
$tcs.Colors(3).RGB  = 6968388   # dk2      44546A
$tcs.Colors(4).RGB  = 15132391  # lt2      E7E6E6
$tcs.Colors(5).RGB  = 13998939  # accent1  5B9BD5
$tcs.Colors(6).RGB  = 3243501   # accent2  ED7D31
$tcs.Colors(7).RGB  = 10855845  # accent3  A5A5A5
$tcs.Colors(8).RGB  = 49407     # accent4  FFC000
$tcs.Colors(9).RGB  = 12874308  # accent5  4472C4
$tcs.Colors(10).RGB = 4697456   # accent6  70AD47
$tcs.Colors(11).RGB = 12673797  # hlink    0563C1
$tcs.Colors(12).RGB = 7491477   # folHlink 954F72
